$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.249.31"
$ws.Range("E2").Value = "  -4.29%  "

# Row 3
$ws.Range("D3").Value = "2.499.57"
$ws.Range("E3").Value = "  -5.47%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.80"
$ws.Range("E5").Value = "  -2.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.88"
$ws.Range("E6").Value = "  -4.82%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -0.97%  "

# Row 9
$ws.Range("D9").Value = "2.496.21"
$ws.Range("E9").Value = "  -5.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -9.60%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("E11").Value = "  -1.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  -4.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  -2.41%  "

# Row 14
$ws.Range("D14").Value = "2.954.00"
$ws.Range("E14").Value = "  -5.47%  "

# Row 15
$ws.Range("D15").Value = "69.188.77"
$ws.Range("E15").Value = "  -4.07%  "

# Row 16
$ws.Range("E16").Value = "  -7.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.67"
$ws.Range("E17").Value = "  -4.77%  "

# Row 18
$ws.Range("D18").Value = "2.499.85"
$ws.Range("E18").Value = "  -4.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.35"
$ws.Range("E19").Value = "  -6.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.75"
$ws.Range("E20").Value = "  -2.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.70"
$ws.Range("E21").Value = "  -6.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("E22").Value = "  -5.69%  "

# Row 23
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  -5.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.05"
$ws.Range("E25").Value = "  -4.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("E26").Value = "  -7.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.86"
$ws.Range("E27").Value = "  -8.06%  "

# Row 28
$ws.Range("E28").Value = "  -4.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.35%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0892"
$ws.Range("E30").Value = "  -6.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  -2.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "466.70"
$ws.Range("E32").Value = "  -6.21%  "

# Row 33
$ws.Range("E33").Value = "  -2.31%  "

# Row 34
$ws.Range("E34").Value = "  -3.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +0.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.91"
$ws.Range("E37").Value = "  -5.89%  "

# Row 38
$ws.Range("E38").Value = "  +0.33%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.33"
$ws.Range("E39").Value = "  -4.89%  "

# Row 40
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.71"
$ws.Range("E41").Value = "  -3.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.314"
$ws.Range("E42").Value = "  -4.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.58"
$ws.Range("E43").Value = "  -9.20%  "

# Row 44
$ws.Range("E44").Value = "  -14.89%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -10.84%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.05"
$ws.Range("E46").Value = "  -2.59%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.31"
$ws.Range("E47").Value = "  -6.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.524"
$ws.Range("E48").Value = "  -4.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.49"
$ws.Range("E49").Value = "  -4.64%  "

# Row 50
$ws.Range("E50").Value = "  -5.46%  "

# Row 51
$ws.Range("E51").Value = "  -2.23%  "
